# Update cryptos list: refresh Price (col D) and Volume(1h) (col E) figures,
# and re-rank a few coins (Maker/ordi swap, THORChain/FraxShare swap) in
# rows 48-51 to reflect the latest GitHub Actions data pull.
# Values that look like plain numbers (e.g. "1.00", "8.38") are written with
# a leading apostrophe so Excel keeps them as text, matching the source data
# which stores prices as text strings (e.g. "42.748.96" style price strings).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '42.748.96'
$ws.Cells.Item(2, 5).Value = '  +0.05%  '
$ws.Cells.Item(3, 4).Value = '2.315.03'
$ws.Cells.Item(3, 5).Value = '  +0.78%  '
$ws.Cells.Item(4, 5).Value = '  -0.09%  '
$ws.Cells.Item(5, 4).Value = '''312.30'
$ws.Cells.Item(5, 5).Value = '  -1.62%  '
$ws.Cells.Item(6, 4).Value = '''106.93'
$ws.Cells.Item(6, 5).Value = '  +2.66%  '
$ws.Cells.Item(7, 4).Value = '''0.622'
$ws.Cells.Item(7, 5).Value = '  -0.90%  '
$ws.Cells.Item(8, 5).Value = '  -0.27%  '
$ws.Cells.Item(9, 5).Value = '  +0.47%  '
$ws.Cells.Item(10, 4).Value = '''40.11'
$ws.Cells.Item(10, 5).Value = '  +1.00%  '
$ws.Cells.Item(11, 4).Value = '''0.0915'
$ws.Cells.Item(11, 5).Value = '  +0.70%  '
$ws.Cells.Item(12, 4).Value = '''8.38'
$ws.Cells.Item(12, 5).Value = '  -1.63%  '
$ws.Cells.Item(13, 4).Value = '''0.108'
$ws.Cells.Item(13, 5).Value = '  -1.28%  '
$ws.Cells.Item(14, 4).Value = '''0.992'
$ws.Cells.Item(14, 5).Value = '  -1.47%  '
$ws.Cells.Item(15, 5).Value = '  -0.28%  '
$ws.Cells.Item(16, 4).Value = '2.664.88'
$ws.Cells.Item(16, 5).Value = '  +0.45%  '
$ws.Cells.Item(17, 4).Value = '2.317.96'
$ws.Cells.Item(17, 5).Value = '  +0.65%  '
$ws.Cells.Item(18, 4).Value = '42.740.29'
$ws.Cells.Item(18, 5).Value = '  +0.10%  '
$ws.Cells.Item(19, 5).Value = '  -0.63%  '
$ws.Cells.Item(20, 5).Value = '  -0.11%  '
$ws.Cells.Item(21, 4).Value = '''13.10'
$ws.Cells.Item(21, 5).Value = '  -12.27%  '
$ws.Cells.Item(22, 4).Value = '''73.60'
$ws.Cells.Item(22, 5).Value = '  -0.46%  '
$ws.Cells.Item(23, 5).Value = '  -1.72%  '
$ws.Cells.Item(24, 4).Value = '''265.71'
$ws.Cells.Item(24, 5).Value = '  -0.29%  '
$ws.Cells.Item(25, 5).Value = '  +1.10%  '
$ws.Cells.Item(26, 4).Value = '''1.01'
$ws.Cells.Item(26, 5).Value = '  +0.07%  '
$ws.Cells.Item(27, 4).Value = '''7.76'
$ws.Cells.Item(27, 5).Value = '  +15.02%  '
$ws.Cells.Item(28, 4).Value = '''11.03'
$ws.Cells.Item(28, 5).Value = '  +0.77%  '
$ws.Cells.Item(29, 4).Value = '''2.26'
$ws.Cells.Item(29, 5).Value = '  -3.57%  '
$ws.Cells.Item(30, 4).Value = '''38.89'
$ws.Cells.Item(30, 5).Value = '  +4.37%  '
$ws.Cells.Item(31, 4).Value = '''22.46'
$ws.Cells.Item(31, 5).Value = '  -0.62%  '
$ws.Cells.Item(32, 4).Value = '''166.24'
$ws.Cells.Item(32, 5).Value = '  +0.24%  '
$ws.Cells.Item(33, 4).Value = '''0.0877'
$ws.Cells.Item(33, 5).Value = '  -0.40%  '
$ws.Cells.Item(34, 5).Value = '  +5.59%  '
$ws.Cells.Item(35, 5).Value = '  -0.98%  '
$ws.Cells.Item(36, 4).Value = '''4.70'
$ws.Cells.Item(36, 5).Value = '  +3.02%  '
$ws.Cells.Item(37, 5).Value = '  -1.34%  '
$ws.Cells.Item(38, 5).Value = '  +1.30%  '
$ws.Cells.Item(39, 4).Value = '''2.85'
$ws.Cells.Item(39, 5).Value = '  +5.57%  '
$ws.Cells.Item(40, 5).Value = '  -1.55%  '
$ws.Cells.Item(41, 4).Value = '''1.62'
$ws.Cells.Item(41, 5).Value = '  +2.34%  '
$ws.Cells.Item(42, 4).Value = '''103.49'
$ws.Cells.Item(42, 5).Value = '  +8.54%  '
$ws.Cells.Item(43, 4).Value = '''71.01'
$ws.Cells.Item(43, 5).Value = '  +0.89%  '
$ws.Cells.Item(44, 5).Value = '  +2.31%  '
$ws.Cells.Item(45, 4).Value = '''12.94'
$ws.Cells.Item(45, 5).Value = '  +5.36%  '
$ws.Cells.Item(46, 4).Value = '''1.00'
$ws.Cells.Item(46, 5).Value = '  -0.43%  '
$ws.Cells.Item(47, 4).Value = '''112.54'
$ws.Cells.Item(47, 5).Value = '  -2.56%  '
$ws.Cells.Item(48, 2).Value = 'ordi'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Cells.Item(48, 4).Value = '''77.07'
$ws.Cells.Item(48, 5).Value = '  -4.80%  '
$ws.Cells.Item(49, 2).Value = 'Maker'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(49, 4).Value = '1.651.99'
$ws.Cells.Item(49, 5).Value = '  -2.51%  '
$ws.Cells.Item(50, 2).Value = 'FraxShare'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(50, 4).Value = '''8.81'
$ws.Cells.Item(50, 5).Value = '  +0.04%  '
$ws.Cells.Item(51, 2).Value = 'THORChain'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Cells.Item(51, 4).Value = '''5.24'
$ws.Cells.Item(51, 5).Value = '  +1.91%  '
